$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 6041
$ws.Range("A2").Value = 5369
$ws.Range("A3").Value = 4508
